# "Elimna EC anteriores y se agregan nuevos, se modifica base de datos"
# The "Periodo Mora" list (rows 16-24, column E) is reversed in order
# (newest period first), and the accompanying "Valor Mora" numeric code
# in column F is re-paired to match the new period ordering.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# New Periodo Mora (col E) values for rows 16..24, newest-first.
$periodos = @("2207", "2206", "2205", "2204", "2203", "2202", "2201", "2112", "2111")

# Matching Valor Mora (col F) numeric codes for rows 16..24.
$valores = @(35112, 35112, 36341, 36341, 36341, 36341, 36341, 36341, 36341)

for ($i = 0; $i -lt $periodos.Length; $i++) {
    $row = 16 + $i
    $ws.Cells.Item($row, 5).Value = $periodos[$i]
    $ws.Cells.Item($row, 6).Value = $valores[$i]
}
